$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 319 (shifts existing rows 319:362 down to 322:365)
$ws.Rows.Item(319).Resize(3).Insert()

# Row 319
$ws.Range("A319").Value = 3
$ws.Range("B319").Value = "Femacal de La Calera"
$ws.Range("C319").Value = "Coquimbo"
$ws.Range("D319").Value = 45173
$ws.Range("E319").Value = 5
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100107
$ws.Range("H319").Value = "Otros"
$ws.Range("I319").Value = 100107002
$ws.Range("J319").Value = "Chirimoya"
$ws.Range("K319").Value = "Cultivar IV Región"
$ws.Range("L319").Value = "Especial"
$ws.Range("M319").Value = 67
$ws.Range("N319").Value = 30000
$ws.Range("O319").Value = 30000
$ws.Range("P319").Value = 30000
$ws.Range("Q319").Value = "$/bandeja 10 kilos"
$ws.Range("R319").Value = "Provincia del Elquí"
$ws.Range("S319").Value = 3000
$ws.Range("T319").Value = 10

# Row 320
$ws.Range("A320").Value = 3
$ws.Range("B320").Value = "Femacal de La Calera"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 45173
$ws.Range("E320").Value = 5
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100107
$ws.Range("H320").Value = "Otros"
$ws.Range("I320").Value = 100107002
$ws.Range("J320").Value = "Chirimoya"
$ws.Range("K320").Value = "Cultivar IV Región"
$ws.Range("L320").Value = "Primera"
$ws.Range("M320").Value = 68
$ws.Range("N320").Value = 30000
$ws.Range("O320").Value = 30000
$ws.Range("P320").Value = 30000
$ws.Range("Q320").Value = "$/bandeja 10 kilos"
$ws.Range("R320").Value = "Provincia del Elquí"
$ws.Range("S320").Value = 3000
$ws.Range("T320").Value = 10

# Row 321
$ws.Range("A321").Value = 3
$ws.Range("B321").Value = "Femacal de La Calera"
$ws.Range("C321").Value = "Coquimbo"
$ws.Range("D321").Value = 45173
$ws.Range("E321").Value = 5
$ws.Range("F321").Value = "Fruta"
$ws.Range("G321").Value = 100107
$ws.Range("H321").Value = "Otros"
$ws.Range("I321").Value = 100107002
$ws.Range("J321").Value = "Chirimoya"
$ws.Range("K321").Value = "Cultivar IV Región"
$ws.Range("L321").Value = "Segunda"
$ws.Range("M321").Value = 56
$ws.Range("N321").Value = 25000
$ws.Range("O321").Value = 25000
$ws.Range("P321").Value = 25000
$ws.Range("Q321").Value = "$/bandeja 10 kilos"
$ws.Range("R321").Value = "Provincia del Elquí"
$ws.Range("S321").Value = 2500
$ws.Range("T321").Value = 10
